$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-23 Friday" "2024-08-24 Saturday"

Replace-Text "42×64=" "49×38="
Replace-Text "33×80=" "72×25="
Replace-Text "63×47=" "25×83="
Replace-Text "53×35=" "32×27="
Replace-Text "94×30=" "11×57="

Replace-Text "66×67=" "26×78="
Replace-Text "98×42=" "48×48="
Replace-Text "51×27=" "18×68="
Replace-Text "33×65=" "15×16="
Replace-Text "29×66=" "79×96="

Replace-Text "29×35=" "21×27="
Replace-Text "58×99=" "75×26="
Replace-Text "35×67=" "30×21="
Replace-Text "45×66=" "86×34="
Replace-Text "52×98=" "88×71="

Replace-Text "11×19=" "99×24="
Replace-Text "55×56=" "28×73="
Replace-Text "42×17=" "47×50="
Replace-Text "89×75=" "91×65="
Replace-Text "32×47=" "92×14="

Replace-Text "80×38=" "18×27="
Replace-Text "16×28=" "71×25="
Replace-Text "21×99=" "69×34="
Replace-Text "93×73=" "16×90="
Replace-Text "25×79=" "31×95="
